$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Move/resize the chart so its anchor spans from the area near C7 down to S37
# (previously anchored roughly J5:U32). Values are in points; the workbook
# uses the default column width (58.4375pt) and default row height (15pt),
# so these translate to the exact target <xdr:from>/<xdr:to> cell+offset
# anchors recorded in the saved drawing XML.
$co = $ws.ChartObjects(1)
$co.Top = 93.37488188976378
$co.Left = 120.24988188976378
$co.Width = 966.1251181102363
$co.Height = 459.3751181102362

# Update the active selection on the sheet from H13 to P4.
$ws.Range("P4").Select()
